$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.349.92"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.35%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.659.79"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.28%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.38%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "220.13"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.506"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.42%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.256"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.91%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.43%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.99"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +4.81%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0848"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.69%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.892.51"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.29%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.659.00"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.99%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.20"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.532"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.10"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.84%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.337.41"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.29%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0735"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "223.25"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +4.85%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.42%  "
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.44"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.78%  "
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = "Chainlink"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.71"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +8.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.45"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +4.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.28"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.39%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.89"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.02%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.73%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +4.05%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.44%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.59%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0513"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.40%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.38"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.16%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.71%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.66%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.262.46"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.61%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0179"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.88%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.85%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.839"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +3.33%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.820"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.37"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.804.55"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.13"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -4.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.90"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.61%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "92.16"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.78%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.61"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.55%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0517"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.10%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.71"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.76%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0980"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.92%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.30%  "
